# Insert a new data row before the current row 740 (2026/12/29, 火, 13, 201):
#   - new row 740: 2026/01/30, 金, 17, 201
#   - all rows from the old 740..781 shift down to 741..782
# This matches the diff: dimension grows from A1:D781 to A1:D782.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 740:781 down by one to make room, then write the new row 740.
$ws.Rows.Item(740).EntireRow.Insert()

# Keep the date column as plain text (matches the rest of the column),
# rather than letting Excel auto-convert it to a date serial value, and
# reset the style back to the sheet's default (no explicit style index),
# matching the look of every other data row.
$ws.Cells.Item(740, 1).NumberFormat = "@"
$ws.Cells.Item(740, 1).Value = "2026/01/30"
$ws.Cells.Item(740, 1).Style = "Normal"
$ws.Cells.Item(740, 2).Value = "金"
$ws.Cells.Item(740, 3).Value = 17
$ws.Cells.Item(740, 4).Value = 201
